$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "phone no with country code"
$ws.Range("A3").Value = "phone no with country code"
$ws.Range("A4").Value = "91xxxxxxxxx"

$ws.Range("C6").Select()
